$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-numeric text cells (Coin names, Links, Volume labels) -- plain assignment keeps them as text.
$ws.Range('B4').Value = 'LEO'
$ws.Range('C4').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('E4').Value = '3LEOLEO'
$ws.Range('B5').Value = 'HuobiToken'
$ws.Range('C5').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E5').Value = '4HuobiTokenHT'
$ws.Range('B6').Value = 'Cronos'
$ws.Range('C6').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E6').Value = '5CronosCRO'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('E7').Value = '6KuCoinTokenKCS'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('E8').Value = '7GateTokenGT'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E9').Value = '8MXTokenMX'
$ws.Range('B10').Value = 'FTXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E10').Value = '9FTXTokenFTT'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('B17').Value = 'BTSEToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('E17').Value = '16BTSETokenBTSE'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E18').Value = '17OneONEBestin24h'
$ws.Range('E27').Value = '26NitroExNTX'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOINWorstin24h'

# Price column (D) holds numeric-looking strings that must stay stored as text,
# matching the original inline-string cells. Force text format first so Excel
# does not auto-convert these into numbers on entry.
$dCells = @('D2','D3','D4','D5','D6','D7','D8','D9','D10','D11','D12','D13','D14','D15','D16','D17','D18','D19','D20','D22','D23','D26','D27','D40','D41','D42','D43','D44','D45','D47')
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '245.26'
$ws.Range('D3').Value = '25.16'
$ws.Range('D4').Value = '3.499'
$ws.Range('D5').Value = '5.044'
$ws.Range('D6').Value = '0.05604'
$ws.Range('D7').Value = '6.555'
$ws.Range('D8').Value = '3.020'
$ws.Range('D9').Value = '0.8141'
$ws.Range('D10').Value = '0.8416'
$ws.Range('D11').Value = '0.1337'
$ws.Range('D12').Value = '0.06956'
$ws.Range('D13').Value = '0.02834'
$ws.Range('D14').Value = '0.09402'
$ws.Range('D15').Value = '0.001511'
$ws.Range('D16').Value = '0.006243'
$ws.Range('D17').Value = '2.091'
$ws.Range('D18').Value = '0.009687'
$ws.Range('D19').Value = '0.3185'
$ws.Range('D20').Value = '0.03253'
$ws.Range('D22').Value = '3.747'
$ws.Range('D23').Value = '0.04680'
$ws.Range('D26').Value = '0.004530'
$ws.Range('D27').Value = '0.00009702'
$ws.Range('D40').Value = '0.03662'
$ws.Range('D41').Value = '0.006231'
$ws.Range('D42').Value = '0.1054'
$ws.Range('D43').Value = '0.002730'
$ws.Range('D44').Value = '0.008157'
$ws.Range('D45').Value = '0.00005298'
$ws.Range('D47').Value = '0.1800'
